$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Not Used"
$ws.Range("C4").Value = "Not Used"
$ws.Range("C5").Value = "ADC[0]"
$ws.Range("C6").Value = "ADC[1]"
$ws.Range("C7").Value = "ADC[2]"

$ws.Range("C7").Select()
